# Rename sheet "Uncut_Sheet" -> "Uncut_Sheet_1" and keep the Print Area
# defined name in sync with the new sheet name. Also update the active
# selection on the sheet to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-point the print area first (while the sheet still has its original
# name) so the underlying defined name is refreshed, then rename the sheet.
$ws.Name = "Uncut_Sheet_1"
$ws.PageSetup.PrintArea = "`$A`$1:`$G`$42"

$ws.Range("B4").Select()
